# Update "想去人数" (column F) figures across sheets, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$ws1.Cells.Item(2, 6).Value = 2692    # F2: 2686 -> 2692
$ws1.Cells.Item(4, 6).Value = 354     # F4: 353 -> 354
$ws1.Cells.Item(5, 6).Value = 1501    # F5: 1498 -> 1501
$ws1.Cells.Item(6, 6).Value = 1139    # F6: 1136 -> 1139
$ws1.Cells.Item(13, 6).Value = 9163   # F13: 9135 -> 9163
$ws1.Cells.Item(14, 6).Value = 396    # F14: 394 -> 396
$ws1.Cells.Item(15, 6).Value = 2499   # F15: 2498 -> 2499
$ws1.Cells.Item(16, 6).Value = 5      # F16: 4 -> 5
$ws1.Cells.Item(20, 6).Value = 630    # F20: 628 -> 630
$ws1.Cells.Item(24, 6).Value = 2088   # F24: 2085 -> 2088
$ws1.Cells.Item(25, 6).Value = 2177   # F25: 2175 -> 2177
$ws1.Cells.Item(27, 6).Value = 1884   # F27: 1882 -> 1884
$ws1.Cells.Item(30, 6).Value = 1533   # F30: 1463 -> 1533
$ws1.Cells.Item(31, 6).Value = 273    # F31: 272 -> 273
$ws1.Cells.Item(32, 6).Value = 154    # F32: 151 -> 154
$ws1.Cells.Item(35, 6).Value = 325    # F35: 324 -> 325
$ws1.Cells.Item(36, 6).Value = 64     # F36: 63 -> 64
$ws1.Cells.Item(37, 6).Value = 291    # F37: 289 -> 291
$ws1.Cells.Item(38, 6).Value = 486    # F38: 487 -> 486
$ws1.Cells.Item(40, 6).Value = 36     # F40: 29 -> 36
$ws1.Cells.Item(41, 6).Value = 468    # F41: 419 -> 468
$ws1.Cells.Item(42, 6).Value = 22     # F42: 14 -> 22
$ws1.Cells.Item(43, 6).Value = 1385   # F43: 1379 -> 1385
$ws1.Cells.Item(44, 6).Value = 297    # F44: 295 -> 297
$ws1.Cells.Item(46, 6).Value = 160    # F46: 8 -> 160
$ws1.Cells.Item(47, 6).Value = 619    # F47: 615 -> 619
$ws1.Cells.Item(49, 6).Value = 297    # F49: 296 -> 297

# --- 演出 (sheet2) ---
$ws2.Cells.Item(7, 6).Value = 19      # F7: 18 -> 19
$ws2.Cells.Item(10, 6).Value = 1      # F10: 0 -> 1

# --- 全部类型 (sheet4) ---
$ws4.Cells.Item(2, 6).Value = 2692    # F2: 2686 -> 2692
$ws4.Cells.Item(3, 6).Value = 354     # F3: 353 -> 354
$ws4.Cells.Item(4, 6).Value = 1501    # F4: 1498 -> 1501
$ws4.Cells.Item(6, 6).Value = 1139    # F6: 1136 -> 1139
$ws4.Cells.Item(10, 6).Value = 9163   # F10: 9135 -> 9163
$ws4.Cells.Item(11, 6).Value = 396    # F11: 394 -> 396
$ws4.Cells.Item(12, 6).Value = 2499   # F12: 2498 -> 2499
$ws4.Cells.Item(14, 6).Value = 5      # F14: 4 -> 5
$ws4.Cells.Item(19, 6).Value = 630    # F19: 628 -> 630
$ws4.Cells.Item(22, 6).Value = 2177   # F22: 2175 -> 2177
$ws4.Cells.Item(23, 6).Value = 1884   # F23: 1882 -> 1884
$ws4.Cells.Item(25, 6).Value = 1533   # F25: 1463 -> 1533
$ws4.Cells.Item(26, 6).Value = 273    # F26: 272 -> 273
$ws4.Cells.Item(27, 6).Value = 154    # F27: 151 -> 154
$ws4.Cells.Item(30, 6).Value = 325    # F30: 324 -> 325
$ws4.Cells.Item(31, 6).Value = 64     # F31: 63 -> 64
$ws4.Cells.Item(32, 6).Value = 291    # F32: 289 -> 291
$ws4.Cells.Item(33, 6).Value = 486    # F33: 487 -> 486
$ws4.Cells.Item(36, 6).Value = 19     # F36: 18 -> 19
$ws4.Cells.Item(38, 6).Value = 36     # F38: 29 -> 36
$ws4.Cells.Item(39, 6).Value = 468    # F39: 419 -> 468
$ws4.Cells.Item(41, 6).Value = 22     # F41: 14 -> 22
$ws4.Cells.Item(42, 6).Value = 1385   # F42: 1379 -> 1385
$ws4.Cells.Item(44, 6).Value = 297    # F44: 295 -> 297
$ws4.Cells.Item(46, 6).Value = 160    # F46: 8 -> 160
$ws4.Cells.Item(47, 6).Value = 619    # F47: 615 -> 619
$ws4.Cells.Item(48, 6).Value = 297    # F48: 296 -> 297
$ws4.Cells.Item(49, 6).Value = 1      # F49: 0 -> 1
